$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ClienteServidor"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "DOM"
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Tables"

$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Layouts"

# ---- Layouts sheet content ----
$ws4.Range("C2").Value = "HEADER"
$ws4.Range("C2:H3").Merge()
$ws4.Range("C2:H3").Interior.Pattern = 1
$ws4.Range("C2:H3").Interior.ThemeColor = 3
$ws4.Range("C2:H3").HorizontalAlignment = -4108

$ws4.Range("C4").Value = "MENU"
$ws4.Range("C4:H4").Merge()
$ws4.Range("C4:H4").Interior.Pattern = 1
$ws4.Range("C4:H4").Interior.ThemeColor = 5
$ws4.Range("C4:H4").HorizontalAlignment = -4108

$ws4.Range("C5").Value = "BARRA DE NAVEGAÇÃO"
$ws4.Range("C5:C15").Merge()
$ws4.Range("C5:C15").Interior.Pattern = 1
$ws4.Range("C5:C15").Interior.ThemeColor = 6
$ws4.Range("C5:C15").HorizontalAlignment = -4108
$ws4.Range("C5:C15").VerticalAlignment = -4108
$ws4.Range("C5:C15").WrapText = $true

$ws4.Range("D5").Value = "CONTEÚDO"
$ws4.Range("D5:G15").Merge()
$ws4.Range("D5:G15").HorizontalAlignment = -4108
$ws4.Range("D5:G15").VerticalAlignment = -4108

$ws4.Range("H5").Value = "ADDS"
$ws4.Range("H5:H15").Merge()
$ws4.Range("H5:H15").Interior.Pattern = 1
$ws4.Range("H5:H15").Interior.ThemeColor = 8
$ws4.Range("H5:H15").HorizontalAlignment = -4108
$ws4.Range("H5:H15").VerticalAlignment = -4108

$ws4.Range("C16").Value = "RODAPÉ"
$ws4.Range("C16:H16").Merge()
$ws4.Range("C16:H16").Interior.Pattern = 1
$ws4.Range("C16:H16").Interior.ThemeColor = 10
$ws4.Range("C16:H16").HorizontalAlignment = -4108

$ws4.Range("H21").Font.Underline = 2

$ws4.PageSetup.PaperSize = 9
$ws4.PageSetup.Orientation = 1

$ws4.Range("C5:H15").Select()
$ws2.Activate()
Write-Host "done"
